$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The four data rows (5-8) get reordered by date (ascending) and their
# row-specific fields updated to match. Columns A,B,C,E,F,G,H,O,R are
# identical across these rows and are left untouched.

# Row 5 -> becomes old Row 7's data (date 2021-04-19 = 44280)
$ws.Range("D5").Value = 44280
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 30
$ws.Range("K5").Value = 25000
$ws.Range("L5").Value = 25000
$ws.Range("M5").Value = 25000
$ws.Range("N5").Value = "$/caja 18 kilos empedrada"
$ws.Range("P5").Value = 1389
$ws.Range("Q5").Value = 18

# Row 6 -> becomes old Row 8's data (date 2021-05-02 = 44293)
$ws.Range("D6").Value = 44293
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 10
$ws.Range("K6").Value = 25000
$ws.Range("L6").Value = 25000
$ws.Range("M6").Value = 25000
$ws.Range("N6").Value = "$/caja 15 kilos empedrada"
$ws.Range("P6").Value = 1667
$ws.Range("Q6").Value = 15

# Row 7 -> becomes old Row 5's data (date 2021-05-24 = 44315)
$ws.Range("D7").Value = 44315
$ws.Range("I7").Value = "Especial"
$ws.Range("J7").Value = 10
$ws.Range("K7").Value = 30000
$ws.Range("L7").Value = 30000
$ws.Range("M7").Value = 30000
$ws.Range("N7").Value = "$/caja 20 kilos empedrada"
$ws.Range("P7").Value = 1500
$ws.Range("Q7").Value = 20

# Row 8 -> becomes old Row 6's data (date 2021-05-24 = 44315)
$ws.Range("D8").Value = 44315
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 20
$ws.Range("K8").Value = 15000
$ws.Range("L8").Value = 15000
$ws.Range("M8").Value = 15000
$ws.Range("N8").Value = "$/caja 15 kilos granel"
$ws.Range("P8").Value = 1000
$ws.Range("Q8").Value = 15
